$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('B2').Value = 'Bitcoin'
$ws.Range('C2').Value = 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc'
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '58.674.07'
$ws.Range('E2').Value = '  -2.19%  '

$ws.Range('B3').Value = 'Ethereum'
$ws.Range('C3').Value = 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.590.25'
$ws.Range('E3').Value = '  -3.16%  '

$ws.Range('B4').Value = 'TetherUSD'
$ws.Range('C4').Value = 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt'
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  -0.02%  '

$ws.Range('B5').Value = 'BNB'
$ws.Range('C5').Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '539.20'
$ws.Range('E5').Value = '  +0.36%  '

$ws.Range('B6').Value = 'Solana'
$ws.Range('C6').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '144.13'
$ws.Range('E6').Value = '  -1.23%  '

$ws.Range('B7').Value = 'USDC'
$ws.Range('C7').Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.999'
$ws.Range('E7').Value = '  -0.05%  '

$ws.Range('B8').Value = 'XRP'
$ws.Range('C8').Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.578'
$ws.Range('E8').Value = '  +0.74%  '

$ws.Range('B9').Value = 'Toncoin'
$ws.Range('C9').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '6.79'
$ws.Range('E9').Value = '  +1.84%  '

$ws.Range('B10').Value = 'Dogecoin'
$ws.Range('C10').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.101'
$ws.Range('E10').Value = '  -3.02%  '

$ws.Range('B11').Value = 'TRON'
$ws.Range('C11').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.139'
$ws.Range('E11').Value = '  +3.45%  '

$ws.Range('B12').Value = 'Cardano'
$ws.Range('C12').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.333'
$ws.Range('E12').Value = '  -1.73%  '

$ws.Range('B13').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C13').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '3.044.96'
$ws.Range('E13').Value = '  -2.85%  '

$ws.Range('B14').Value = 'WrappedBTC'
$ws.Range('C14').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '58.578.47'
$ws.Range('E14').Value = '  -2.23%  '

$ws.Range('B15').Value = 'Avalanche'
$ws.Range('C15').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '20.78'
$ws.Range('E15').Value = '  -2.63%  '

$ws.Range('B16').Value = 'WrappedEther'
$ws.Range('C16').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.594.26'
$ws.Range('E16').Value = '  -2.66%  '

$ws.Range('B17').Value = 'ShibaInu'
$ws.Range('C17').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.0000134'
$ws.Range('E17').Value = '  -1.43%  '

$ws.Range('B18').Value = 'Polkadot'
$ws.Range('C18').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '4.48'
$ws.Range('E18').Value = '  +1.09%  '

$ws.Range('B19').Value = 'BitcoinCash'
$ws.Range('C19').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '337.05'
$ws.Range('E19').Value = '  -2.43%  '

$ws.Range('B20').Value = 'Chainlink'
$ws.Range('C20').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '10.10'
$ws.Range('E20').Value = '  -1.93%  '

$ws.Range('B21').Value = 'Uniswap'
$ws.Range('C21').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.18'
$ws.Range('E21').Value = '  -3.27%  '

$ws.Range('B22').Value = 'Dai'
$ws.Range('C22').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.999'
$ws.Range('E22').Value = '  +0.05%  '

$ws.Range('B23').Value = 'Litecoin'
$ws.Range('C23').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '67.04'
$ws.Range('E23').Value = '  -0.82%  '

$ws.Range('B24').Value = 'Polygon'
$ws.Range('C24').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.420'
$ws.Range('E24').Value = '  +0.76%  '

$ws.Range('B25').Value = 'Binance-PegBSC-USD'
$ws.Range('C25').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.996'
$ws.Range('E25').Value = '  -0.15%  '

$ws.Range('B26').Value = 'Kaspa'
$ws.Range('C26').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.159'
$ws.Range('E26').Value = '  -4.62%  '

$ws.Range('B27').Value = 'InternetComputer(DFINITY)'
$ws.Range('C27').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.09'
$ws.Range('E27').Value = '  -3.53%  '

$ws.Range('B28').Value = 'PEPE'
$ws.Range('C28').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.0₃0738'
$ws.Range('E28').Value = '  -2.35%  '

$ws.Range('B29').Value = 'USDe'
$ws.Range('C29').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.999'
$ws.Range('E29').Value = '  -0.03%  '

$ws.Range('B30').Value = 'PancakeSwap'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.65'
$ws.Range('E30').Value = '  -1.09%  '

$ws.Range('B31').Value = 'Monero'
$ws.Range('C31').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '156.62'
$ws.Range('E31').Value = '  +4.06%  '

$ws.Range('B32').Value = 'Aptos'
$ws.Range('C32').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '5.95'
$ws.Range('E32').Value = '  +0.32%  '

$ws.Range('B33').Value = 'EthereumClassic'
$ws.Range('C33').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '18.87'
$ws.Range('E33').Value = '  -1.48%  '

$ws.Range('B34').Value = 'NEARProtocol'
$ws.Range('C34').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.92'
$ws.Range('E34').Value = '  -3.44%  '

$ws.Range('B35').Value = 'OKB'
$ws.Range('C35').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '37.01'
$ws.Range('E35').Value = '  -0.79%  '

$ws.Range('B36').Value = 'ImmutableX'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.10'
$ws.Range('E36').Value = '  -5.12%  '

$ws.Range('B37').Value = 'SuiNetwork'
$ws.Range('C37').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.849'
$ws.Range('E37').Value = '  +2.61%  '

$ws.Range('B38').Value = 'Fetch.AI'
$ws.Range('C38').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.826'
$ws.Range('E38').Value = '  -2.22%  '

$ws.Range('B39').Value = 'Stacks'
$ws.Range('C39').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.42'
$ws.Range('E39').Value = '  -3.67%  '

$ws.Range('B40').Value = 'Filecoin'
$ws.Range('C40').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.61'
$ws.Range('E40').Value = '  -0.38%  '

$ws.Range('B41').Value = 'Bittensor'
$ws.Range('C41').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '283.06'
$ws.Range('E41').Value = '  -2.98%  '

$ws.Range('B42').Value = 'FirstDigitalUSD'
$ws.Range('C42').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.998'
$ws.Range('E42').Value = '  -0.22%  '

$ws.Range('B43').Value = 'Mantle'
$ws.Range('C43').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.593'
$ws.Range('E43').Value = '  -2.36%  '

$ws.Range('B44').Value = 'Hedera'
$ws.Range('C44').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0535'
$ws.Range('E44').Value = '  -0.75%  '

$ws.Range('B45').Value = 'WhiteBITCoin'
$ws.Range('C45').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '10.63'
$ws.Range('E45').Value = '  -1.05%  '

$ws.Range('B46').Value = 'Stellar'
$ws.Range('C46').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0942'
$ws.Range('E46').Value = '  -1.57%  '

$ws.Range('B47').Value = 'EnergySwap'
$ws.Range('C47').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '18.61'
$ws.Range('E47').Value = '  -4.39%  '

$ws.Range('B48').Value = 'VeChain'
$ws.Range('C48').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0228'
$ws.Range('E48').Value = '  +0.45%  '

$ws.Range('B49').Value = 'Maker'
$ws.Range('C49').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.920.58'
$ws.Range('E49').Value = '  -2.99%  '

$ws.Range('B50').Value = 'InjectiveProtocol'
$ws.Range('C50').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '18.02'
$ws.Range('E50').Value = '  -2.88%  '

$ws.Range('B51').Value = 'RenderToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '4.44'
$ws.Range('E51').Value = '  -3.64%  '
